$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "SavingErrors" row (currently row 11) values before the row above it is removed.
$ws.Range("B11").Value = 22
$ws.Range("C11").Value = 24

# Delete the "ImportNativeFiles" row (row 5) data, shifting only columns A:E up by
# one (the F:I summary block stays anchored to its rows).
$ws.Range("A5:E5").Delete(-4162)

# Update the selected cell shown in the sheet view.
$ws.Range("D10").Select()
